$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

# Row 21 (existing): B21 currently holds the placeholder "string" value;
# replace it with the new "a borrar" label.
$ws.Range("B21").Value = "a borrar"

# Row 22: new user row. "2222" must be stored as *text*, not auto-converted
# to a number by Excel's smart typing, so it is written through a
# Text-formatted helper cell + copy/paste-values (then the helper is wiped).
$ws.Range("A22").Value = 12313124
$ws.Range("C22").Value = "string"
$ws.Range("D22").Value = "string"

# Row 23: duplicate of row 22
$ws.Range("A23").Value = 12313124
$ws.Range("C23").Value = "string"
$ws.Range("D23").Value = "string"

$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value = "2222"
$helper.Copy()
$ws.Range("B22").PasteSpecial(-4163)
$ws.Range("B23").PasteSpecial(-4163)
$helper.Clear()

# Row 24: final new user row, with an empty LibrosPrestados value.
# A no-op formatting touch (re-asserting the existing default) materialises
# the otherwise-empty D24 cell without tagging it with a distinct style.
$ws.Range("A24").Value = 99999991
$ws.Range("B24").Value = "solito despues de borrar"
$ws.Range("C24").Value = "string"
$ws.Range("D24").Font.Bold = $false
